$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.662.70"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "1.903.31"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "'245.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.44%  "
$ws.Range("E6").Value = "  +2.38%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "'42.73"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("D9").Value = "'0.338"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.68%  "
$ws.Range("D10").Value = "'0.0708"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "2.179.63"
$ws.Range("E12").Value = "  +3.07%  "
$ws.Range("D13").Value = "'12.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.19%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.692"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.889.23"
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("D16").Value = "'4.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("D17").Value = "35.631.75"
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").Value = "'72.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.06%  "
$ws.Range("D19").Value = "0.0₃0813"
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("D20").Value = "'245.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").Value = "'4.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.48%  "
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").Value = "'2.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.48%  "
$ws.Range("D25").Value = "'171.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").Value = "'2.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +30.25%  "
$ws.Range("E27").Value = "  +7.22%  "
$ws.Range("D28").Value = "'18.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("D30").Value = "'0.958"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +27.48%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.64%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.0567"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.36%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("D34").Value = "'4.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.70%  "
$ws.Range("D35").Value = "'1.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.14%  "
$ws.Range("D36").Value = "'2.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.48%  "
$ws.Range("E37").Value = "  +8.50%  "
$ws.Range("E38").Value = "  +3.61%  "
$ws.Range("D39").Value = "'0.0206"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'91.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("D41").Value = "1.366.44"
$ws.Range("E41").Value = "  +1.23%  "
$ws.Range("D42").Value = "'15.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.51%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.0597"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.31%  "
$ws.Range("B44").Value = "Gas"
$ws.Range("C44").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D44").Value = "'13.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +46.30%  "
$ws.Range("D45").Value = "'2.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.30%  "
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").Value = "'6.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.18%  "
$ws.Range("E48").Value = "  +35.41%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").Value = "2.093.28"
$ws.Range("E50").Value = "  +2.92%  "
$ws.Range("D51").Value = "'3.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.37%  "
